# Apply updated state-rankings data: prison guard/corrections officer
# contribution groups removed from the underlying candidate/party data.
# total_dollars (B) drops for affected states; dollars_per1000ppl (D) is
# recomputed as B/C*1000; rows are re-ranked/re-sorted by D descending
# (E), and F reflects the population rank for the resulting row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$finalData = @(
    @("NY", 14923643.8, 19618453, 760.6942198755428, 1, 4),
    @("CA", 26431158.66, 39148760, 675.1467648017459, 2, 1),
    @("AK", 412361.63, 738516, 558.3651945252371, 3, 47),
    @("NJ", 4339350.65, 8881845, 488.5641046426728, 4, 11),
    @("FL", 6278758.97, 20598139, 304.8216622870639, 5, 3),
    @("RI", 317855.0, 1056611, 300.82499614333, 6, 43),
    @("NV", 832020.29, 2922849, 284.6607163079585, 7, 34),
    @("WA", 2025249.12, 7294336, 277.646809798726, 8, 13),
    @("DE", 238219.0, 949495, 250.8902100590314, 9, 45),
    @("VA", 2101776.64, 8413774, 249.8018891403549, 10, 12),
    @("MA", 1396997.86, 6830193, 204.5327064696415, 11, 15),
    @("PA", 2579650.08, 12791181, 201.6741128125699, 12, 6),
    @("MN", 942407.51, 5527358, 170.49872832554, 13, 22),
    @("MD", 1017934.63, 6003435, 169.5586993113109, 14, 19),
    @("LA", 770424.2, 4663616, 165.1988928762574, 15, 25),
    @("IL", 1896057.56, 12821497, 147.8811374366036, 16, 5),
    @("MI", 1415030.0, 9957488, 142.1071258132573, 17, 10),
    @("TX", 3834138.02, 27885195, 137.4972640499735, 18, 2),
    @("OH", 1066361.97, 11641879, 91.59706693395455, 19, 7),
    @("HI", 122385.94, 1422029, 86.06430670541882, 20, 40),
    @("WI", 491294.0, 5778394, 85.02258586036189, 21, 20),
    @("OK", 323350.0, 3918137, 82.52646602198953, 22, 28),
    @("NH", 100345.0, 1343622, 74.68246277598908, 23, 41),
    @("CT", 243441.19, 3581504, 67.97177666142493, 24, 29),
    @("MO", 340045.0, 6090062, 55.83604895976429, 25, 18),
    @("AZ", 383107.0, 6946685, 55.14961452836857, 26, 14),
    @("AL", 266775.44, 4864680, 54.83925766956922, 27, 24),
    @("NE", 103557.51, 1904760, 54.36774711774712, 28, 37),
    @("MT", 49789.95, 1041732, 47.79535427537984, 29, 44),
    @("IN", 292786.13, 6637426, 44.11139649617186, 30, 17),
    @("OR", 173880.0, 4081943, 42.59736110964803, 31, 27),
    @("KS", 97000.0, 2908776, 33.34735985170394, 32, 35),
    @("TN", 219233.76, 6651089, 32.96208485557779, 33, 16),
    @("WV", 57850.0, 1829054, 31.62837182499806, 34, 38),
    @("NC", 262756.19, 10155624, 25.87297343816589, 35, 9),
    @("IA", 77435.0, 3132499, 24.71988019788673, 36, 30),
    @("ID", 28675.0, 1687809, 16.98948162973417, 37, 39),
    @("UT", 50402.48, 3045350, 16.5506362158701, 38, 31),
    @("GA", 146475.19, 10297484, 14.22436684533814, 39, 8),
    @("AR", 39326.44, 2990671, 13.14970453119049, 40, 32),
    @("CO", 69601.62, 5531141, 12.58359170377324, 41, 21),
    @("MS", 24033.1, 2988762, 8.041155501843239, 42, 33),
    @("SC", 29186.76, 4955925, 5.889265878720925, 43, 23),
    @("ME", 7800.0, 1332813, 5.852283853773935, 44, 42),
    @("NM", 5800.0, 2092434, 2.7718914909622, 45, 36),
    @("VT", 750.0, 624977, 1.200044161625148, 46, 48),
    @("ND", 600.0, 752201, 0.7976591363212758, 47, 46),
    @("KY", 1000.0, 4440204, 0.2252148775146367, 48, 26),
)

$rowIndex = 2
foreach ($rec in $finalData) {
    $ws.Cells.Item($rowIndex, 1).Value2 = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value2 = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value2 = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value2 = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value2 = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value2 = $rec[5]
    $rowIndex = $rowIndex + 1
}

Write-Output "Updated $($finalData.Count) state rows"
